# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# existing header style used by the other header cells (B1:G1) and
# filling the data rows (H2:H8) with 0, matching column F/G's numeric
# (unstyled) cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header style (bold, bordered, centered) from the adjacent
# "sum" header cell (G1) onto the new "Save" header cell (H1), then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new data column with 0 for every existing data row.
$ws.Range("H2:H8").Value = 0
